$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(139).Insert()

$ws.Range("A139").Value = 10
$ws.Range("B139").Value = 'Vega Modelo de Temuco'
$ws.Range("C139").Value = 'La Araucanía'
$ws.Range("D139").Value = 45033
$ws.Range("E139").Value = 9
$ws.Range("F139").Value = 'Fruta'
$ws.Range("G139").Value = 100104
$ws.Range("H139").Value = 'Frutos de pepita'
$ws.Range("I139").Value = 100104003
$ws.Range("J139").Value = 'Membrillo'
$ws.Range("K139").Value = 'Champion'
$ws.Range("L139").Value = 'Primera'
$ws.Range("M139").Value = 120
$ws.Range("N139").Value = 14000
$ws.Range("O139").Value = 14000
$ws.Range("P139").Value = 14000
$ws.Range("Q139").Value = '$/bandeja 18 kilos granel'
$ws.Range("R139").Value = 'Región de O''Higgins'
$ws.Range("S139").Value = 778
$ws.Range("T139").Value = 18
